$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell address -> new text value (all values are text/labels in this sheet)
$updates = [ordered]@{
    "D2" = "68.608.22"
    "E2" = "  -1.56%  "
    "D3" = "3.859.09"
    "E3" = "  -0.99%  "
    "E4" = "  +0.01%  "
    "D5" = "602.93"
    "E5" = "  -0.25%  "
    "D6" = "169.49"
    "E6" = "  -0.92%  "
    "D7" = "3.858.54"
    "E7" = "  -1.01%  "
    "E9" = "  -1.35%  "
    "E10" = "  -1.76%  "
    "E11" = "  +1.46%  "
    "E12" = "  -2.13%  "
    "E13" = "  +4.81%  "
    "D14" = "37.13"
    "D15" = "4.502.84"
    "E15" = "  -1.05%  "
    "D16" = "3.853.32"
    "E16" = "  -1.14%  "
    "D17" = "68.780.56"
    "E17" = "  -1.29%  "
    "E18" = "  -1.04%  "
    "E19" = "  -2.98%  "
    "E20" = "  -0.85%  "
    "E21" = "  +0.88%  "
    "D22" = "471.09"
    "E22" = "  -4.02%  "
    "E23" = "  -1.63%  "
    "E24" = "  -1.44%  "
    "D25" = "83.56"
    "E25" = "  -2.09%  "
    "E26" = "  -2.82%  "
    "E27" = "  -1.36%  "
    "D28" = "10.18"
    "E28" = "  +0.29%  "
    "E29" = "  +0.17%  "
    "E30" = "  -0.39%  "
    "D31" = "4.009.35"
    "E32" = "  -2.34%  "
    "D33" = "31.54"
    "E33" = "  -1.45%  "
    "E34" = "  -4.19%  "
    "E35" = "  -3.14%  "
    "D36" = "3.822.93"
    "E36" = "  -1.08%  "
    "E37" = "  -2.05%  "
    "D38" = "3.75"
    "E38" = "  +9.47%  "
    "E39" = "  -1.49%  "
    "E40" = "  -1.99%  "
    "D41" = "5.95"
    "E41" = "  -2.96%  "
    "E42" = "  +0.06%  "
    "E43" = "  -4.02%  "
    "D44" = "1.99"
    "E44" = "  -5.27%  "
    "D45" = "8.74"
    "E45" = "  +0.56%  "
    "D46" = "0.000296"
    "E46" = "  +7.23%  "
    "D47" = "417.95"
    "E47" = "  -3.87%  "
    "E48" = "  -0.01%  "
    "E49" = "  -1.94%  "
    "E50" = "  -1.48%  "
    "D51" = "142.01"
    "E51" = "  -0.71%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking strings (e.g. "37.13") are not
    # auto-converted to numbers by Excel, preserving the original text cell type.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Reset the cell style back to the workbook default so no incidental
    # formatting change is introduced.
    $cell.Style = "Normal"
}
